# protocol.xlsx - "Updated StatsWriter to OpenCL 1.2"
#
# The author appended one new diary entry (row 65 on Sheet1: a date in
# column A and a note in column B) and left the cell cursor on B70
# afterwards. This reproduces that content edit through the Excel
# object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet   # Sheet1 is tabSelected="1" -> the active sheet

# New diary row: date 2013-04-09 (serial 41373) + the new shared-string note.
$ws.Range("A65").Value = 41373
$ws.Range("B65").Value = "Updated device informations queried in StatsWriter to OpenCL 1.2 (not tested)"

# Author's cursor ended up on B70 (was B68 before the edit).
$ws.Range("B70").Select()
